# Generate Report for Handback
# Updates the handoff/handback timestamps for the
# 45b7ba3a-cb1e-46bf-9291-4edc66d2a960 file (row 2 on each sheet)
# following a fresh handback report generation.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the 45b7ba3a file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-30 09:06:25"

# zh-cn sheet: Correspond Handoff Datetime (H2) / Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-30 09:06:21"
$wsZhCn.Range("K2").Value = "2016-08-30 09:06:38"

# de-de sheet: Correspond Handoff Datetime (H2) / Correspond Handback DateTime (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-30 09:06:25"
$wsDeDe.Range("K2").Value = "2016-08-30 09:06:45"
